{"js": "// 1) Typo fix: \"get it's own\" -> \"get its own\" in the Entertainment-law\n//    paragraph (FirstParagraph style, \"What is Entertainment Law?\" section).\nconst typoResults = context.document.body.search(\"get it's own\", { matchCase: true });\ntypoResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < typoResults.items.length; i++) {\n  typoResults.items[i].insertText(\"get its own\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Re-style the \"--Frank Zappa\" paragraph from \"First Paragraph\" to\n//    \"Block Text\" so it matches the other quote attributions above it.\nconst zappaResults = context.document.body.search(\"Frank Zappa\", { matchCase: true });\nzappaResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < zappaResults.items.length; i++) {\n  const para = zappaResults.items[i].paragraphs.getFirst();\n  para.style = \"Block Text\";\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Typo fix: \"get it's own\" -> \"get its own\" in the Entertainment-law\n#    paragraph (FirstParagraph style, \"What is Entertainment Law?\" section).\n$find = $d.Content.Find\n$find.Text = \"get it's own\"\n$find.Replacement.Text = \"get its own\"\n$find.Execute(\n    $find.Text,        # FindText\n    $false,            # MatchCase\n    $false,            # MatchWholeWord\n    $false,            # MatchWildcards\n    $false,            # MatchSoundsLike\n    $false,            # MatchAllWordForms\n    $true,             # Forward\n    1,                 # Wrap (wdFindContinue)\n    $false,            # Format\n    $find.Replacement.Text,  # ReplaceWith\n    2                  # Replace (wdReplaceAll)\n)\n\n# 2) Re-style the \"--Frank Zappa\" paragraph from \"First Paragraph\" to\n#    \"Block Text\" so it matches the other quote attributions above it.\n$rng = $d.Content\n$find2 = $rng.Find\n$find2.Text = \"Frank Zappa\"\n$found = $find2.Execute()\nif ($found) {\n    $para = $rng.Paragraphs(1)\n    $para.Style = \"Block Text\"\n}\n"}
